# Auto-generated edit script: Add data for 2023-11-04
# Updates column J (2023 totals) across Citywide Totals, By Neighborhood, and
# individual neighborhood sheets to reflect one additional day of violent crime data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 6503
$ws.Range("J3").Value = 6897
$ws.Range("H4").Value = 1706
$ws.Range("J4").Value = 1494
$ws.Range("J6").Value = 9128
$ws.Range("H7").Value = 26017
$ws.Range("J7").Value = 24552

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 422
$ws.Range("J6").Value = 541
$ws.Range("J7").Value = 1547

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J3").Value = 181
$ws.Range("J7").Value = 485

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 366
$ws.Range("J6").Value = 384
$ws.Range("J7").Value = 1103

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 128
$ws.Range("J6").Value = 79
$ws.Range("J7").Value = 357

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 220
$ws.Range("J7").Value = 752

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 180
$ws.Range("J3").Value = 171
$ws.Range("J7").Value = 611

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 149
$ws.Range("J7").Value = 376

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 195
$ws.Range("J4").Value = 113
$ws.Range("J5").Value = 75
$ws.Range("J7").Value = 717
$ws.Range("J8").Value = 1547
$ws.Range("J12").Value = 51
$ws.Range("J14").Value = 130
$ws.Range("J16").Value = 100
$ws.Range("J17").Value = 34
$ws.Range("J18").Value = 210
$ws.Range("J19").Value = 719
$ws.Range("J20").Value = 510
$ws.Range("J22").Value = 58
$ws.Range("I25").Value = 131
$ws.Range("J29").Value = 1336
$ws.Range("J33").Value = 1103
$ws.Range("J34").Value = 113
$ws.Range("J36").Value = 334
$ws.Range("J37").Value = 752
$ws.Range("J41").Value = 166
$ws.Range("J42").Value = 1054
$ws.Range("J43").Value = 210
$ws.Range("J44").Value = 186
$ws.Range("J49").Value = 157
$ws.Range("J50").Value = 148
$ws.Range("J52").Value = 620
$ws.Range("J54").Value = 465
$ws.Range("J58").Value = 15
$ws.Range("J60").Value = 143
$ws.Range("J61").Value = 27
$ws.Range("H63").Value = 264
$ws.Range("I63").Value = 251
$ws.Range("J63").Value = 85
$ws.Range("J64").Value = 164
$ws.Range("J65").Value = 611
$ws.Range("J66").Value = 74
$ws.Range("J67").Value = 928
$ws.Range("J71").Value = 79
$ws.Range("J73").Value = 236
$ws.Range("J75").Value = 73
$ws.Range("J77").Value = 178
$ws.Range("J78").Value = 292
$ws.Range("J79").Value = 690
$ws.Range("J83").Value = 485
$ws.Range("J84").Value = 204
$ws.Range("J85").Value = 1022
$ws.Range("J86").Value = 159
$ws.Range("J88").Value = 252
$ws.Range("J89").Value = 319
$ws.Range("J90").Value = 260
$ws.Range("J92").Value = 78
$ws.Range("J94").Value = 260
$ws.Range("J95").Value = 357
$ws.Range("J97").Value = 219
$ws.Range("J99").Value = 376
$ws.Range("H101").Value = 26017
$ws.Range("J101").Value = 24552

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 345
$ws.Range("J7").Value = 928

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J2").Value = 62
$ws.Range("J7").Value = 204

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J2").Value = 28
$ws.Range("J3").Value = 30
$ws.Range("J7").Value = 157

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 218
$ws.Range("J7").Value = 465

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 407
$ws.Range("J4").Value = 71
$ws.Range("J6").Value = 336
$ws.Range("J7").Value = 1336

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 169
$ws.Range("J7").Value = 719

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 186

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J2").Value = 34
$ws.Range("J6").Value = 97
$ws.Range("J7").Value = 166

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 206
$ws.Range("J6").Value = 561
$ws.Range("J7").Value = 1054

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J2").Value = 78
$ws.Range("J7").Value = 292

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 233
$ws.Range("J6").Value = 204
$ws.Range("J7").Value = 690

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J6").Value = 140
$ws.Range("J7").Value = 510

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J2").Value = 57
$ws.Range("J7").Value = 210

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 107
$ws.Range("J7").Value = 334

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 222
$ws.Range("J6").Value = 232
$ws.Range("J7").Value = 717

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J3").Value = 28
$ws.Range("J4").Value = 7
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 113

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J2").Value = 50
$ws.Range("J3").Value = 48
$ws.Range("J4").Value = 19
$ws.Range("J7").Value = 260

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 131

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J2").Value = 39
$ws.Range("J7").Value = 148

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 74

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 76
$ws.Range("J7").Value = 236

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J4").Value = 15
$ws.Range("J7").Value = 195

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J2").Value = 37
$ws.Range("J3").Value = 21
$ws.Range("J7").Value = 219

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J2").Value = 26
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J6").Value = 122
$ws.Range("J7").Value = 252

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J4").Value = 33
$ws.Range("J7").Value = 319

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J4").Value = 87
$ws.Range("J7").Value = 159

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 92
$ws.Range("J7").Value = 260

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J6").Value = 40
$ws.Range("J7").Value = 143

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J2").Value = 24
$ws.Range("J6").Value = 125
$ws.Range("J7").Value = 210

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J4").Value = 68
$ws.Range("J6").Value = 296
$ws.Range("J7").Value = 1022

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("J2").Value = 28
$ws.Range("J7").Value = 58

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("J2").Value = 22
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J4").Value = 23
$ws.Range("J6").Value = 264
$ws.Range("J7").Value = 620

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J6").Value = 49
$ws.Range("J7").Value = 113

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("J2").Value = 9
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("J3").Value = 2
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 100

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 15

Write-Host "Applied 187 cell updates across 54 sheets."
